{"js": "// Replace the date line and each division-problem cell's text.\n// Each 'old' string is unique in the document at the moment it is searched\n// (searches run strictly in document order), so a plain text search +\n// Replace is unambiguous even though some 'new' values coincide with\n// other cells' 'old' values elsewhere in the table.\nconst replacements = [\n  [\"2023-11-24 Friday\", \"2023-11-25 Saturday\"],\n  [\"71\u00f76=11, 5\", \"13\u00f75=2, 3\"],\n  [\"18\u00f73=6, 0\", \"16\u00f79=1, 7\"],\n  [\"34\u00f79=3, 7\", \"66\u00f77=9, 3\"],\n  [\"30\u00f73=10, 0\", \"64\u00f78=8, 0\"],\n  [\"72\u00f73=24, 0\", \"89\u00f76=14, 5\"],\n  [\"59\u00f78=7, 3\", \"54\u00f79=6, 0\"],\n  [\"50\u00f79=5, 5\", \"52\u00f76=8, 4\"],\n  [\"41\u00f74=10, 1\", \"27\u00f74=6, 3\"],\n  [\"16\u00f76=2, 4\", \"51\u00f77=7, 2\"],\n  [\"10\u00f74=2, 2\", \"30\u00f73=10, 0\"],\n  [\"48\u00f77=6, 6\", \"16\u00f78=2, 0\"],\n  [\"52\u00f72=26, 0\", \"78\u00f74=19, 2\"],\n  [\"87\u00f72=43, 1\", \"89\u00f79=9, 8\"],\n  [\"15\u00f72=7, 1\", \"29\u00f74=7, 1\"],\n  [\"75\u00f78=9, 3\", \"32\u00f74=8, 0\"],\n  [\"65\u00f75=13, 0\", \"25\u00f74=6, 1\"],\n  [\"30\u00f72=15, 0\", \"25\u00f74=6, 1\"],\n  [\"61\u00f73=20, 1\", \"86\u00f77=12, 2\"],\n  [\"81\u00f74=20, 1\", \"58\u00f72=29, 0\"],\n  [\"78\u00f78=9, 6\", \"15\u00f77=2, 1\"],\n  [\"60\u00f73=20, 0\", \"53\u00f73=17, 2\"],\n  [\"24\u00f77=3, 3\", \"31\u00f79=3, 4\"],\n  [\"50\u00f73=16, 2\", \"53\u00f79=5, 8\"],\n  [\"34\u00f73=11, 1\", \"90\u00f73=30, 0\"],\n  [\"70\u00f76=11, 4\", \"21\u00f77=3, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: \"${oldText}\"`);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each division-problem cell's text using\n# Word's Find/Replace (Range.Find.Execute), one pair at a time.\n#\n# Each 'old' string is unique in the document at the moment it is\n# searched (pairs are applied strictly in document order, matching the\n# order the table's cells appear in), so a plain Find/Replace is\n# unambiguous even though some 'new' values coincide with other cells'\n# 'old' values elsewhere in the table (e.g. '25\u00f74=6, 1' is produced\n# twice, and '30\u00f73=10, 0' is both an original value and a later result).\n\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$pairs = @(\n  @(\"2023-11-24 Friday\", \"2023-11-25 Saturday\"),\n  @(\"71\u00f76=11, 5\", \"13\u00f75=2, 3\"),\n  @(\"18\u00f73=6, 0\", \"16\u00f79=1, 7\"),\n  @(\"34\u00f79=3, 7\", \"66\u00f77=9, 3\"),\n  @(\"30\u00f73=10, 0\", \"64\u00f78=8, 0\"),\n  @(\"72\u00f73=24, 0\", \"89\u00f76=14, 5\"),\n  @(\"59\u00f78=7, 3\", \"54\u00f79=6, 0\"),\n  @(\"50\u00f79=5, 5\", \"52\u00f76=8, 4\"),\n  @(\"41\u00f74=10, 1\", \"27\u00f74=6, 3\"),\n  @(\"16\u00f76=2, 4\", \"51\u00f77=7, 2\"),\n  @(\"10\u00f74=2, 2\", \"30\u00f73=10, 0\"),\n  @(\"48\u00f77=6, 6\", \"16\u00f78=2, 0\"),\n  @(\"52\u00f72=26, 0\", \"78\u00f74=19, 2\"),\n  @(\"87\u00f72=43, 1\", \"89\u00f79=9, 8\"),\n  @(\"15\u00f72=7, 1\", \"29\u00f74=7, 1\"),\n  @(\"75\u00f78=9, 3\", \"32\u00f74=8, 0\"),\n  @(\"65\u00f75=13, 0\", \"25\u00f74=6, 1\"),\n  @(\"30\u00f72=15, 0\", \"25\u00f74=6, 1\"),\n  @(\"61\u00f73=20, 1\", \"86\u00f77=12, 2\"),\n  @(\"81\u00f74=20, 1\", \"58\u00f72=29, 0\"),\n  @(\"78\u00f78=9, 6\", \"15\u00f77=2, 1\"),\n  @(\"60\u00f73=20, 0\", \"53\u00f73=17, 2\"),\n  @(\"24\u00f77=3, 3\", \"31\u00f79=3, 4\"),\n  @(\"50\u00f73=16, 2\", \"53\u00f79=5, 8\"),\n  @(\"34\u00f73=11, 1\", \"90\u00f73=30, 0\"),\n  @(\"70\u00f76=11, 4\", \"21\u00f77=3, 0\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceOne)\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n"}
